$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial for each data row.
# Every row whose value is the old serial 45205 moves forward one day to 45206.
$usedRange = $ws.UsedRange
$lastRow = $usedRange.Row + $usedRange.Rows.Count - 1
if ($lastRow -lt 100) { $lastRow = 100 }

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Range("C$row")
    if ($cell.Value2 -eq 45205) {
        $cell.Value = 45206
    }
}
